# Update database: shift each "twelve months ended" column left by one
# fiscal year and append the newly reported 1401/12 year, per the
# commit "update database and change read_price algorithm".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Header labels (rows 8 and 24): shift fiscal-year labels forward ---
$headerRow1 = 8
$headerRow2 = 24
$headers = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)
$cols = @("E", "F", "G", "H", "I")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + $headerRow1).Value = $headers[$i]
    $ws.Range($cols[$i] + $headerRow2).Value = $headers[$i]
}

# --- Data rows: shift each row's 5 yearly figures left, append new value ---
# Row => new [E, F, G, H, I] values (post-shift, newest year appended at I)
$rowValues = @{
    10 = @(132648, 153327, 514249, 695206, 782031)
    12 = @(0, 0, 0, 0, 0)
    13 = @(1028, 2139, 1635, 2581, 10148)
    14 = @(751, 1630, 1046, 1780, 5971)
    15 = @(1637, 1335, 1389, 3842, 4385)
    16 = @(1601, 1046, 2697, 3505, 5799)
    17 = @(51767, 65060, 88024, 132156, 233763)
    19 = @(61784, 129646, 182910, 254216, 317718)
    20 = @(251216, 354183, 791950, 1093286, 1359815)
    26 = @(122, 110, 107, 107, 108)
    27 = @(642, 632, 637, 640, 730)
}

foreach ($row in $rowValues.Keys) {
    $vals = $rowValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}
